$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Create the new "TableCaption" paragraph style (based on Normal1)
#    so table captions ("{{table:...}}" placeholders) look captiony.
# ------------------------------------------------------------------
$tc = $d.Styles.Add("TableCaption", 1)
$tc.BaseStyle = "Normal1"
$tc.QuickStyle = $true

$tc.Font.NameAscii = "Calibri"
$tc.Font.NameFarEast = "Calibri"
$tc.Font.NameBi = "Calibri"
$tc.Font.Name = "Calibri"
$tc.Font.Bold = $true
$tc.Font.Color = 3355443
$tc.Font.Size = 9

$tc.ParagraphFormat.LineSpacingRule = 0

# ------------------------------------------------------------------
# 2. Find the paragraph that holds the {{table:ecosystems}} merge
#    field and re-style it with the new TableCaption style, dropping
#    all of its direct paragraph/run formatting.
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text -like "*{{table:ecosystems}}*") {
        $target = $para
        break
    }
}

$target.Style = "TableCaption"

# ------------------------------------------------------------------
# 3. Re-insert the "_GoBack" bookmark at the start of that paragraph.
#    Word only ever keeps a single "_GoBack" bookmark, so adding the
#    new one automatically removes the old one that used to sit in
#    the {{table:protection}} paragraph, and every other bookmark's
#    w:id is shifted up to make room.
# ------------------------------------------------------------------
$goBackRange = $target.Range.Duplicate
$goBackRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $goBackRange)
